$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0682643427741467
$ws.Range("C2").Value = 0.0442992011619463
$ws.Range("D2").Value = 0.239651416122004
$ws.Range("E2").Value = 0.00798838053740015
$ws.Range("F2").Value = 0.9840232389252
$ws.Range("G2").Value = 0.0217864923747277
$ws.Range("H2").Value = 0.758169934640523
$ws.Range("I2").Value = 0.37763253449528
$ws.Range("J2").Value = 0.989832970225127
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0.957879448075526
$ws.Range("M2").Value = 0.948438634713144
$ws.Range("N2").Value = 0.0319535221496006
$ws.Range("O2").Value = 0.0348583877995643
$ws.Range("P2").Value = 0.10239651416122
$ws.Range("Q2").Value = 0.0312273057371097
$ws.Range("R2").Value = 0.0341321713870733
$ws.Range("S2").Value = 0.872912127814089
$ws.Range("T2").Value = 0.226579520697168
$ws.Range("U2").Value = 0.38562091503268
$ws.Range("V2").Value = 0.100217864923747
$ws.Range("W2").Value = 0.00798838053740015
$ws.Range("X2").Value = 0.0159767610748003

$ws.Range("B3").Value = 0.793028322440087
$ws.Range("C3").Value = 0.904865649963689
$ws.Range("D3").Value = 0.0530137981118373
$ws.Range("E3").Value = 0.037763253449528
$ws.Range("F3").Value = 0.00217864923747277
$ws.Range("G3").Value = 0.00145243282498184
$ws.Range("H3").Value = 0.00217864923747277
$ws.Range("I3").Value = 0.000726216412490922
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.00145243282498184
$ws.Range("L3").Value = 0.00580973129992738
$ws.Range("M3").Value = 0.0261437908496732
$ws.Range("N3").Value = 0.00798838053740015
$ws.Range("O3").Value = 0.0130718954248366
$ws.Range("P3").Value = 0.0348583877995643
$ws.Range("Q3").Value = 0.135802469135802
$ws.Range("R3").Value = 0.0341321713870733
$ws.Range("S3").Value = 0.0646332607116921
$ws.Range("T3").Value = 0.146695715323166
$ws.Range("U3").Value = 0.0305010893246187
$ws.Range("V3").Value = 0.0392156862745098
$ws.Range("W3").Value = 0.0137981118373275
$ws.Range("X3").Value = 0.0363108206245461

$ws.Range("B4").Value = 0.0203340595497458
$ws.Range("C4").Value = 0.0137981118373275
$ws.Range("D4").Value = 0.498910675381264
$ws.Range("E4").Value = 0.950617283950617
$ws.Range("F4").Value = 0.00944081336238199
$ws.Range("G4").Value = 0.973129992737836
$ws.Range("H4").Value = 0.237472766884532
$ws.Range("I4").Value = 0.612926652142338
$ws.Range("J4").Value = 0.00726216412490922
$ws.Range("K4").Value = 0.998547567175018
$ws.Range("L4").Value = 0.0341321713870733
$ws.Range("M4").Value = 0.0152505446623094
$ws.Range("N4").Value = 0.959331880900508
$ws.Range("O4").Value = 0.0355846042120552
$ws.Range("P4").Value = 0.856935366739288
$ws.Range("Q4").Value = 0.037763253449528
$ws.Range("R4").Value = 0.00798838053740015
$ws.Range("S4").Value = 0.0588235294117647
$ws.Range("T4").Value = 0.581699346405229
$ws.Range("U4").Value = 0.0617283950617284
$ws.Range("V4").Value = 0.848220769789397
$ws.Range("W4").Value = 0.046477850399419
$ws.Range("X4").Value = 0.916485112563544

$ws.Range("B5").Value = 0.11837327523602
$ws.Range("C5").Value = 0.037037037037037
$ws.Range("D5").Value = 0.207697893972404
$ws.Range("E5").Value = 0.00363108206245461
$ws.Range("F5").Value = 0.00435729847494553
$ws.Range("G5").Value = 0.00363108206245461
$ws.Range("H5").Value = 0.00217864923747277
$ws.Range("I5").Value = 0.00871459694989107
$ws.Range("J5").Value = 0.00290486564996369
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0.00217864923747277
$ws.Range("M5").Value = 0.0101670297748729
$ws.Range("N5").Value = 0.000726216412490922
$ws.Range("O5").Value = 0.916485112563544
$ws.Range("P5").Value = 0.00580973129992738
$ws.Range("Q5").Value = 0.793754538852578
$ws.Range("R5").Value = 0.923747276688453
$ws.Range("S5").Value = 0.00217864923747277
$ws.Range("T5").Value = 0.0442992011619463
$ws.Range("U5").Value = 0.522149600580973
$ws.Range("V5").Value = 0.0123456790123457
$ws.Range("W5").Value = 0.931735657225853
$ws.Range("X5").Value = 0.0312273057371097
